# Apply updated cryptocurrency market data to the active worksheet.
# Values that Excel would otherwise auto-parse as numbers (plain decimals in
# the "Price" column) are entered with a leading apostrophe to force them to
# stay as text (matching the original inline-string cells), and the resulting
# quote-prefix cell style is then reset to "Normal" so formatting is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "38.287.18"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +3.09%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "2.123.61"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +3.57%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "235.83"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.84%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "0.627"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +1.07%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "58.49"; ForceText = $true }
    @{ Cell = "E7"; Value = "  +2.46%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.395"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +3.43%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.0783"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +3.74%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.104"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +2.08%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "2.439.46"; ForceText = $false }
    @{ Cell = "E12"; Value = "  +3.73%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "14.78"; ForceText = $true }
    @{ Cell = "E13"; Value = "  +3.70%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "21.86"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +5.21%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "0.794"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +3.10%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "5.28"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +2.70%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "2.119.94"; ForceText = $false }
    @{ Cell = "E17"; Value = "  +3.41%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "38.247.47"; ForceText = $false }
    @{ Cell = "E18"; Value = "  +3.12%  "; ForceText = $false }
    @{ Cell = "E19"; Value = "  -1.62%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "70.78"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +2.92%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "0.0₃0830"; ForceText = $false }
    @{ Cell = "E21"; Value = "  +2.88%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "229.63"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +2.36%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -0.11%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "2.42"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -0.44%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.42"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +2.17%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "168.69"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +1.56%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "0.143"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +13.85%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "9.08"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +3.81%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  -0.83%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "19.60"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +3.32%  "; ForceText = $false }
    @{ Cell = "E31"; Value = "  +1.77%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "4.68"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +5.40%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "2.63"; ForceText = $true }
    @{ Cell = "E33"; Value = "  +4.42%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "0.0631"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +2.88%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "4.64"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +1.82%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "3.52"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +8.06%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  +4.78%  "; ForceText = $false }
    @{ Cell = "E38"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "5.52"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -3.06%  "; ForceText = $false }
    @{ Cell = "E40"; Value = "  +8.76%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "2.95"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +0.23%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "97.76"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +1.22%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "0.0216"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +3.34%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "1.465.34"; ForceText = $false }
    @{ Cell = "E44"; Value = "  -0.56%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "1.17"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +0.38%  "; ForceText = $false }
    @{ Cell = "B46"; Value = "InjectiveProtocol"; ForceText = $false }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false }
    @{ Cell = "D46"; Value = "16.23"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +7.42%  "; ForceText = $false }
    @{ Cell = "B47"; Value = "ARBITRUM"; ForceText = $false }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; ForceText = $false }
    @{ Cell = "D47"; Value = "1.08"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +6.18%  "; ForceText = $false }
    @{ Cell = "B48"; Value = "FTXToken"; ForceText = $false }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; ForceText = $false }
    @{ Cell = "D48"; Value = "4.14"; ForceText = $true }
    @{ Cell = "E48"; Value = "  -6.76%  "; ForceText = $false }
    @{ Cell = "E49"; Value = "  +4.04%  "; ForceText = $false }
    @{ Cell = "E50"; Value = "  +2.32%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "2.322.49"; ForceText = $false }
    @{ Cell = "E51"; Value = "  +3.65%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
